# Daily attendance processing - 2025-11-08 04:47:26
# Normalize the "Recorded By" (column G) entries so that any "System" /
# "system" token in the comma-separated list of recorders is moved from
# the end to the front of the list (the whole token order is reversed).
# Rows whose "Recorded By" value does not include a "System" token (e.g.
# a lone email, or an email/email pair) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value -eq "") { continue }

    $parts = $value -split ", "

    if ($parts.Count -le 1) { continue }

    $hasSystem = $false
    foreach ($part in $parts) {
        if ($part.Trim().ToLower() -eq "system") {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $reversed = $parts[($parts.Count - 1)..0]
        $newValue = $reversed -join ", "
        $cell.Value = $newValue
    }
}
